# Season-record columns: Wins / Losses / Ties
# The scraper previously only pulled team statistics; this adds the
# season W-L-T record as three new trailing columns (AD:AF).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header formatting of the existing header row (bold, centered,
# bordered) by copying the format from the adjacent "Unnamed: 28" header.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Every player row on the roster shares the same team season record.
$lastRow = 55
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 57
    $ws.Cells.Item($r, 31).Value = 105
    $ws.Cells.Item($r, 32).Value = 0
}
